$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Phase 1: append all new paragraphs with their plain text, no
# italics yet. Doing every InsertAfter up front (before any
# Font.Italic toggling) avoids the "type-after-italic-text
# inherits italic" behaviour bleeding into later paragraphs.
# We remember the absolute [start,end) span of every run that
# must end up italic, and italicize them all in Phase 2 once all
# the text exists.
# ---------------------------------------------------------------
$italicSpans = New-Object System.Collections.ArrayList

# Anchor: the last paragraph currently in the main body
# (the BILAGA 1 title paragraph) -- new paragraphs get appended
# after it in order, each becoming the new anchor for the next.
$anchor = $d.Paragraphs.Last.Range

# --- paragraph 0: style='Heading1' ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Heading1'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Knärot – ekologi samt krav på livsmiljön')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 1: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 2: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Samuel Johnsons doktorsavhandling ')
$r.Collapse(0)
$r.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$r.Collapse(0)
$r.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter('Vidare ')
$r.Collapse(0)
$r.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 3: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ')
$r.Collapse(0)
$r.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 4: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 5: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 6: style='Heading2' ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Heading2'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Referenser - knärot')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 7: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('de Graaf M & Roberts M.R., 2009. ')
$r.Collapse(0)
$r.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter('Forest Ecology and Management 257, 1014-1025')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 8: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ')
$r.Collapse(0)
$r.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter('Ecological Applications, 22, 2049-2064 ')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 9: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ')
$r.Collapse(0)
$r.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter('Landscape Ecology, 37, sid 1839-1853')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 10: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ')
$r.Collapse(0)
$r.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter('Journal of Applied Ecology. 51, 53-62.')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 11: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('Skogsstyrelsen, 2022. ')
$r.Collapse(0)
$r.InsertAfter('Vägledning för hänsyn till knärot. ')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$r.Collapse(0)
$anchor = $p.Range

# --- paragraph 12: style=None ---
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = 'Normal'
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter('SLU Artdatabanken, 2021. ')
$r.Collapse(0)
$r.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$italicSpans.Add(@($r.Start, $r.End)) | Out-Null
$r.Collapse(0)
$r.InsertAfter('SLU Artdatabanken, Uppsala ')
$r.Collapse(0)
$anchor = $p.Range

# ---------------------------------------------------------------
# Phase 2: italicize the remembered spans now that all text in
# the new section has been typed.
# ---------------------------------------------------------------
foreach ($span in $italicSpans) {
    $sub = $d.Range($span[0], $span[1])
    $sub.Font.Italic = 1
}


# --- Fix the date stamp in the first-page header (header3.xml / rId12) ---
$sec = $d.Sections.First
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "2023-09-15", 2) | Out-Null
    }
}
